$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1 and give it the same formatting as the other
# header cells (e.g. G1: bold, bordered, centered) by copying the format
# from the neighboring header cell.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the "Save" column data values for the existing data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
